$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the data of rows 15, 16 and 18 (row 17 is left alone):
#   new row15 <- old row16
#   new row16 <- old row18
#   new row18 <- old row15
# The "Publik kommentar" note (column AC) travels along with the row it
# belongs to, so it moves from row 15 to row 18.
#
# Note: this COM shim's plain `.Value` getter does not reliably return the
# underlying value (it can hand back the property's own reflection info
# instead), so reads go through `.Value2` - except for column I, whose
# numeric-looking content is stored as text in the workbook and must be
# read back with `.Text` to avoid silently turning it into a number.

$cols = @("A","B","D","E","F","G","H","I","P","Q","R")

function Get-RowValues($rowNum) {
    $vals = @{}
    foreach ($c in $cols) {
        $cell = $ws.Range("${c}${rowNum}")
        if ($c -eq "I") {
            $vals[$c] = $cell.Text
        }
        else {
            $vals[$c] = $cell.Value2
        }
    }
    $vals["AC"] = $ws.Range("AC${rowNum}").Value2
    return $vals
}

$old15 = Get-RowValues 15
$old16 = Get-RowValues 16
$old18 = Get-RowValues 18

function Set-RowValues($rowNum, $vals) {
    foreach ($c in $cols) {
        $cell = $ws.Range("${c}${rowNum}")
        if ($c -eq "I") {
            # Column I stores numeric-looking values as text in the source file.
            $cell.NumberFormat = "@"
        }
        $cell.Value = $vals[$c]
    }
}

Set-RowValues 15 $old16
Set-RowValues 16 $old18
Set-RowValues 18 $old15

$ws.Range("AC15").Value = ""
$ws.Range("AC16").Value = ""
if ($old15["AC"]) {
    $ws.Range("AC18").Value = $old15["AC"]
}
